$d = $word.ActiveDocument

function Find-ParagraphIndex($doc, $substr1, $substr2) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $t = $doc.Paragraphs.Item($i).Range.Text
        if ($t -like $substr1 -and $t -like $substr2) {
            return $i
        }
    }
    return -1
}

# --- Edit 1: split the run ending in "...for each subgoal. " so "subgoal" is
#     wrapped in spell-check proofErr markers, matching the author's correction. ---
$idx1 = Find-ParagraphIndex $d "*subgoal*" "*Continue with this*"
if ($idx1 -lt 0) { throw "Could not find target paragraph for edit 1 (subgoal)" }
$p1 = $d.Paragraphs.Item($idx1)
$p1.Range.InsertXML(@'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="3E6B0310" w14:textId="2DC3C432" w:rsidR="00554479" w:rsidRDefault="00554479" w:rsidP="00554479"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="240"/><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">a) </w:t></w:r><w:r w:rsidR="00A27B7B"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">Determine the pattern </w:t></w:r><w:r w:rsidR="004F6CFF"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">so as to know where it restarts, allowing one to calculate where exactly her counting ends in larger increments.  In this case, 1-10 end on the first finger, then 11-20 end on the ring finger, then 21-30 end on the ring finger, and then 31-40 end on the first finger before counting restarts at the thumb, ending again at the first finger at 50.  Thus every 50 counted digits ends on the first finger.  Continue with this to determine the finger landed upon for each </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>subgoal</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@)

# --- Edit 2: rewrite problem 3's solution section (questions 4b through the
#     trailing empty paragraph) to add the worked-out answer (5a) and the
#     test-case write-up (5b), moving the _GoBack bookmark to the new end. ---
$idx2 = Find-ParagraphIndex $d "*Does each solution meet the goals*" "*Will each solution work for ALL cases*"
if ($idx2 -lt 0) { throw "Could not find target paragraph for edit 2 (evaluate solutions)" }
$pStart = $d.Paragraphs.Item($idx2)
$pEnd = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$r.InsertXML(@'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="2A6152FB" w14:textId="70C6CBFE" w:rsidR="00554479" w:rsidRDefault="00554479" w:rsidP="00554479"><w:pPr><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="220"/><w:tab w:val="left" w:pos="720"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="240"/><w:ind w:hanging="720"/><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">a) Does each solution meet the goals? </w:t></w:r><w:r w:rsidR="000062F5"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> Yes. </w:t></w:r><w:r w:rsidR="000062F5"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/></w:rPr><w:t> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">b) Will each solution work for ALL cases? </w:t></w:r><w:r w:rsidR="000062F5"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> Yes. </w:t></w:r></w:p><w:p><w:pPr><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="220"/><w:tab w:val="left" w:pos="720"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="240"/><w:ind w:hanging="720"/><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">5)  Choose a solution and develop a plan to implement it. </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/></w:rPr><w:t> </w:t></w:r></w:p><w:p><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="240"/><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">a) </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">Determine the pattern so as to know where it restarts, allowing one to calculate where exactly her counting ends in larger increments.  In this case, 1-10 end on the first finger, then 11-20 end on the ring finger, then 21-30 end on the ring finger, and then 31-40 end on the first finger before counting restarts at the thumb, ending again at the first finger at 50.  Thus every 50 counted digits ends on the first finger. </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> Using this to determine the answer, ending on 10 leaves the girl on the first finger. Ending on 100 also leaves the girl ending on the first finger. Following this all the way up to 1000, you also see her ending on her first finger.  </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:br/></w:r></w:p><w:p><w:pPr><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="220"/><w:tab w:val="left" w:pos="720"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="240"/><w:ind w:hanging="720"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">b) Describe some test cases you tried out to make sure it works. </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:br/><w:t xml:space="preserve">To test this, I counted on my fingers up until 50 to determine the pattern and to see where the counting “resets,” thus determining if the results were repeatable. As this was verified and the math works, it should work every time. </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@)

Write-Output "done: idx1=$idx1 idx2=$idx2 finalCount=$($d.Paragraphs.Count)"
